# Change isPartOfBook property to isPartOfBookChapter (#63)
# Updates the row describing the "isPartOfBook" property so that it now
# describes "isPartOfBookChapter" instead (name, labels, comments and the
# linked object type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "isPartOfBookChapter"
$ws.Range("B26").Value = "Part of chapter"
$ws.Range("C26").Value = "Teil des Kapitels"
$ws.Range("D26").Value = "Appartient au chapitre"
$ws.Range("E26").Value = "Appartiene al capitolo"
$ws.Range("G26").Value = "Belongs to following chapter"
$ws.Range("H26").Value = "Gehört zu folgendem Kapitel"
$ws.Range("I26").Value = "Appartient au chapitre suivant"
$ws.Range("J26").Value = "Appartiene al seguente capitolo"
$ws.Range("M26").Value = ":BookChapter"

[void]$ws.Range("A26").Select()
